$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "41.485.91"
Set-TextValue "E2" "  -1.86%  "

# Row 3
Set-TextValue "D3" "2.196.24"
Set-TextValue "E3" "  -1.55%  "

# Row 4
Set-TextValue "E4" "  -0.09%  "

# Row 5
Set-TextValue "D5" "239.64"
Set-TextValue "E5" "  -1.74%  "

# Row 6
Set-TextValue "D6" "0.620"
Set-TextValue "E6" "  -1.59%  "

# Row 7
Set-TextValue "D7" "71.97"
Set-TextValue "E7" "  -2.92%  "

# Row 9
Set-TextValue "D9" "0.589"
Set-TextValue "E9" "  -4.42%  "

# Row 10
Set-TextValue "D10" "41.43"
Set-TextValue "E10" "  -4.09%  "

# Row 11
Set-TextValue "D11" "0.0943"
Set-TextValue "E11" "  -2.71%  "

# Row 12
Set-TextValue "B12" "Polkadot"
Set-TextValue "C12" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "6.88"
Set-TextValue "E12" "  -4.23%  "

# Row 13
Set-TextValue "B13" "TRON"
Set-TextValue "C13" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D13" "0.102"
Set-TextValue "E13" "  -1.20%  "

# Row 14
Set-TextValue "D14" "2.525.40"
Set-TextValue "E14" "  -1.65%  "

# Row 15
Set-TextValue "D15" "14.16"
Set-TextValue "E15" "  -1.32%  "

# Row 16
Set-TextValue "D16" "0.826"
Set-TextValue "E16" "  -3.12%  "

# Row 17
Set-TextValue "D17" "2.192.23"
Set-TextValue "E17" "  -2.49%  "

# Row 18
Set-TextValue "D18" "41.349.46"
Set-TextValue "E18" "  -1.93%  "

# Row 19
Set-TextValue "E19" "  -7.81%  "

# Row 20
Set-TextValue "D20" "6.10"
Set-TextValue "E20" "  -1.38%  "

# Row 21
Set-TextValue "D21" "71.40"
Set-TextValue "E21" "  -1.12%  "

# Row 22
Set-TextValue "D22" "10.36"
Set-TextValue "E22" "  +2.31%  "

# Row 23
Set-TextValue "D23" "228.45"
Set-TextValue "E23" "  -1.17%  "

# Row 24
Set-TextValue "D24" "2.01"
Set-TextValue "E24" "  -8.17%  "

# Row 25
Set-TextValue "E25" "  +0.07%  "

# Row 26
Set-TextValue "D26" "11.15"
Set-TextValue "E26" "  -6.15%  "

# Row 27
Set-TextValue "E27" "  +0.55%  "

# Row 28
Set-TextValue "E28" "  -2.58%  "

# Row 29
Set-TextValue "E29" "  -0.71%  "

# Row 30
Set-TextValue "D30" "166.46"
Set-TextValue "E30" "  -0.28%  "

# Row 31
Set-TextValue "D31" "20.28"
Set-TextValue "E31" "  -2.43%  "

# Row 32
Set-TextValue "B32" "InjectiveProtocol"
Set-TextValue "C32" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D32" "31.84"
Set-TextValue "E32" "  +7.18%  "

# Row 33
Set-TextValue "B33" "Hedera"
Set-TextValue "C33" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D33" "0.0794"
Set-TextValue "E33" "  -1.39%  "

# Row 34
Set-TextValue "D34" "5.27"
Set-TextValue "E34" "  -6.88%  "

# Row 35
Set-TextValue "D35" "0.122"
Set-TextValue "E35" "  -2.38%  "

# Row 36
Set-TextValue "E36" "  -9.78%  "

# Row 37
Set-TextValue "D37" "4.19"
Set-TextValue "E37" "  -5.04%  "

# Row 38
Set-TextValue "E38" "  -3.10%  "

# Row 39
Set-TextValue "D39" "12.88"
Set-TextValue "E39" "  -2.74%  "

# Row 40
Set-TextValue "D40" "2.10"
Set-TextValue "E40" "  -3.13%  "

# Row 41
Set-TextValue "B41" "THORChain"
Set-TextValue "C41" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D41" "5.53"
Set-TextValue "E41" "  -1.64%  "

# Row 42
Set-TextValue "B42" "MultiversX"
Set-TextValue "C42" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D42" "61.57"
Set-TextValue "E42" "  -3.43%  "

# Row 43
Set-TextValue "E43" "  -3.72%  "

# Row 44
Set-TextValue "E44" "  -3.08%  "

# Row 45
Set-TextValue "B45" "Aave"
Set-TextValue "C45" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D45" "100.87"
Set-TextValue "E45" "  -4.20%  "

# Row 46
Set-TextValue "B46" "Cronos"
Set-TextValue "C46" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D46" "0.0987"
Set-TextValue "E46" "  -2.98%  "

# Row 47
Set-TextValue "E47" "  -1.39%  "

# Row 48
Set-TextValue "D48" "1.16"
Set-TextValue "E48" "  -1.48%  "

# Row 49
Set-TextValue "E49" "  -5.93%  "

# Row 50
Set-TextValue "E50" "  -1.62%  "

# Row 51
Set-TextValue "D51" "2.402.39"
Set-TextValue "E51" "  -1.62%  "
